# Daily attendance processing - 2025-11-26 11:47:44
#
# Normalizes the "Recorded By" column (G) on the "Session Analysis Results"
# sheet: wherever the recorder list is exactly "System, dnasr281@gmail.com"
# or "admin@admin.com, dnasr281@gmail.com", move dnasr281@gmail.com to the
# front of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Text

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    } elseif ($val -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, admin@admin.com"
    }
}
